$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row (row 36) describing source J0856+0224
$ws.Range("A36").Value = "J0856+0224"
$ws.Range("B36").Value = 5.55
$ws.Range("B36").HorizontalAlignment = -4152   # xlRight, matches existing column B style
$ws.Range("C36").Value = 899.82
$ws.Range("D36").Value = -1.18
$ws.Range("E36").Value = "Drouart+20"

# Update the view: scroll so row 22 is at top and select D31,
# matching the saved worksheet state from the edit.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D31").Select() | Out-Null
